$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: column B renamed from "Host Number In Wave" to "Host Name In Wave" ---
$ws.Range("B1").Value = "Host Name In Wave"

# --- Row 2 (Second Flow 1 / Capture / Linux src1) ---
$ws.Range("B2").Value = "psp-MyLinSecondFlow-src1"
$ws.Range("D2").Value = "psp-MyLin1-src-IMAGE"

# --- Row 3 (Second Flow 1 / Capture / Linux src2) ---
$ws.Range("B3").Value = "psp-MyLinSecondFlow-src2"
$ws.Range("D3").Value = "psp-MyLin2-src-IMAGE"

# --- Row 4 (Second Flow / Existing System / Windows src1) ---
$ws.Range("B4").Value = "psp-MyWinSecondFlow-src1"
$ws.Range("G4").Value = "172.29.30.156"

# --- Row 5 (Second Flow / Existing System / Windows src2) ---
$ws.Range("B5").Value = "psp-MyWinSecondFlow-src2"
$ws.Range("G5").Value = "172.29.30.219"
$ws.Range("G5").Font.Color = 2236962

# --- Row 6 (Second Flow 1 / Existing System / Linux src2) ---
$ws.Range("B6").Value = "psp-MyLinSecondFlow-src2"
$ws.Range("G6").Value = "172.29.30.218"

# --- Row 7 (Second Flow 1 / Existing System / Linux src1) ---
$ws.Range("B7").Value = "psp-MyLinSecondFlow-src1"
$ws.Range("G7").Value = "172.29.30.155"

# --- Column widths: widen column B for the longer host names, add column I width ---
$ws.Columns.Item(2).ColumnWidth = 23.44140625
$ws.Columns.Item(9).ColumnWidth = 10.21875

# --- Selection: previously selected I6:I7, now just I7 ---
$ws.Range("I7").Select()

"done"
